# Update the "想去人数" (interest count) figures in column F across sheets,
# as output was regenerated at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Sheet "展览"
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 771
$ws1.Range("F5").Value = 237
$ws1.Range("F7").Value = 171
$ws1.Range("F8").Value = 6339
$ws1.Range("F10").Value = 410
$ws1.Range("F12").Value = 5433
$ws1.Range("F15").Value = 1204
$ws1.Range("F24").Value = 3892

# Sheet "演出"
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F2").Value = 97

# Sheet "全部类型"
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 97
$ws4.Range("F4").Value = 771
$ws4.Range("F6").Value = 237
$ws4.Range("F8").Value = 171
$ws4.Range("F9").Value = 6339
$ws4.Range("F11").Value = 410
$ws4.Range("F13").Value = 5433
$ws4.Range("F16").Value = 1204
$ws4.Range("F25").Value = 3892
